$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "Product Name"
$ws.Range("B4").Value = "Samsung 138 cm (55 Inches) Super 6 Series 4K UHD LED Smart TV UA55NU6100 (Black) (2019 model)"

$ws.Range("B12").Select()
